# Scheduled runner: refresh Universalis market-price snapshots (currentAveragePrice /
# NQ / HQ, Leve buy/profit columns) for the affected Leve rows across each crafting job sheet.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets("ALC")
$ws.Range("H28").Value = 924.4666999999999
$ws.Range("I28").Value = 967.2222
$ws.Range("J28").Value = 860.3333
$ws.Range("K28").Value = 967.2222
$ws.Range("L28").Value = 860.3333
$ws.Range("M28").Value = -482.2222
$ws.Range("N28").Value = -1830.3333
$ws.Range("H58").Value = 652.4286
$ws.Range("I58").Value = 591.6667
$ws.Range("K58").Value = 1775.0001
$ws.Range("M58").Value = -1625.0001
$ws.Range("H116").Value = 4277.4287
$ws.Range("I116").Value = 4073.8333
$ws.Range("K116").Value = 4073.8333
$ws.Range("M116").Value = -631.8332999999998
$ws.Range("H131").Value = 1007
$ws.Range("I131").Value = 1007
$ws.Range("K131").Value = 3021
$ws.Range("M131").Value = 2019
$ws.Range("H132").Value = 1527.8695
$ws.Range("I132").Value = 1370.0454
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4110.1362
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1580.1362
$ws.Range("N132").Value = -20060
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets("ARM")
$ws.Range("H14").Value = 50000
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets("BSM")
$ws.Range("H20").Value = 912.5
$ws.Range("I20").Value = 883.5
$ws.Range("K20").Value = 883.5
$ws.Range("M20").Value = -636.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets("CRP")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H31").Value = 1998
$ws.Range("I31").Value = 1998
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1998
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1703
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1998
$ws.Range("I34").Value = 1998
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1998
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1796
$ws.Range("N34").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets("CUL")
$ws.Range("H2").Value = 41.35
$ws.Range("J2").Value = 86
$ws.Range("L2").Value = 516
$ws.Range("N2").Value = -742
$ws.Range("I4").Value = 490.92307
$ws.Range("J4").Value = 253968.25
$ws.Range("K4").Value = 1472.76921
$ws.Range("L4").Value = 761904.75
$ws.Range("M4").Value = -1360.76921
$ws.Range("N4").Value = -762128.75
$ws.Range("H6").Value = 196.75
$ws.Range("I6").Value = 196.75
$ws.Range("K6").Value = 590.25
$ws.Range("M6").Value = -477.25
$ws.Range("H22").Value = 65375.5
$ws.Range("I22").Value = 83833.336
$ws.Range("J22").Value = 10002
$ws.Range("K22").Value = 251500.008
$ws.Range("L22").Value = 30006
$ws.Range("M22").Value = -251331.008
$ws.Range("N22").Value = -30344
$ws.Range("H24").Value = 500
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 1500
$ws.Range("N24").Value = -1960
$ws.Range("M24").ClearContents()
$ws.Range("H27").Value = 65375.5
$ws.Range("I27").Value = 83833.336
$ws.Range("J27").Value = 10002
$ws.Range("K27").Value = 251500.008
$ws.Range("L27").Value = 30006
$ws.Range("M27").Value = -251398.008
$ws.Range("N27").Value = -30210

# --- GSM ---
$ws = $wb.Worksheets("GSM")
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 10000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 10000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -9471
$ws.Range("N25").ClearContents()
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -754
$ws.Range("H107").Value = 8719.833000000001
$ws.Range("J107").Value = 16933
$ws.Range("L107").Value = 16933
$ws.Range("N107").Value = -20773
$ws.Range("H113").Value = 2513.889
$ws.Range("I113").Value = 1776.5714
$ws.Range("J113").Value = 5094.5
$ws.Range("K113").Value = 1776.5714
$ws.Range("L113").Value = 5094.5
$ws.Range("M113").Value = 393.4286
$ws.Range("N113").Value = -9434.5

# --- LTW ---
$ws = $wb.Worksheets("LTW")
$ws.Range("H16").Value = 130.5
$ws.Range("I16").Value = 130.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 130.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 39.5
$ws.Range("N16").ClearContents()
$ws.Range("H38").Value = 32999
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H50").Value = 10000
$ws.Range("I50").Value = 10000
$ws.Range("K50").Value = 10000
$ws.Range("M50").Value = -9363
$ws.Range("H61").Value = 1750.5
$ws.Range("I61").Value = 1538.125
$ws.Range("K61").Value = 1538.125
$ws.Range("M61").Value = -1336.125
$ws.Range("H100").Value = 10781.333
$ws.Range("I100").Value = 10781.333
$ws.Range("K100").Value = 10781.333
$ws.Range("M100").Value = -10240.333
$ws.Range("H113").Value = 1750.5
$ws.Range("I113").Value = 1538.125
$ws.Range("K113").Value = 1538.125
$ws.Range("M113").Value = 631.875
$ws.Range("H131").Value = 39999
$ws.Range("J131").Value = 39999
$ws.Range("L131").Value = 39999
$ws.Range("N131").Value = -50079

# --- WVR ---
$ws = $wb.Worksheets("WVR")
$ws.Range("H9").Value = 5003
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 6
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 134
$ws.Range("N9").Value = -10280
$ws.Range("H113").Value = 460.2
$ws.Range("I113").Value = 476
$ws.Range("J113").Value = 444.4
$ws.Range("K113").Value = 1428
$ws.Range("L113").Value = 1333.2
$ws.Range("M113").Value = 742
$ws.Range("N113").Value = -5673.2
